$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 72/73, pushing existing data (old rows 72-184)
# down to become rows 74-186.
$ws.Rows("72:73").Insert()

# New row 72: Packham's Triumph / Primera, 03-12-2021, volumen 700
$ws.Cells.Item(72, 1).Value = 4
$ws.Cells.Item(72, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(72, 3).Value = "Los Lagos"
$ws.Cells.Item(72, 4).Value = 44533
$ws.Cells.Item(72, 5).Value = 10
$ws.Cells.Item(72, 6).Value = "Fruta"
$ws.Cells.Item(72, 7).Value = 100104
$ws.Cells.Item(72, 8).Value = "Frutos de pepita"
$ws.Cells.Item(72, 9).Value = 100104005
$ws.Cells.Item(72, 10).Value = "Pera"
$ws.Cells.Item(72, 11).Value = "Packham's Triumph"
$ws.Cells.Item(72, 12).Value = "Primera"
$ws.Cells.Item(72, 13).Value = 700
$ws.Cells.Item(72, 14).Value = 15000
$ws.Cells.Item(72, 15).Value = 16000
$ws.Cells.Item(72, 16).Value = 15500
$ws.Cells.Item(72, 17).Value = "`$/caja 15 kilos empedrada"
$ws.Cells.Item(72, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(72, 19).Value = 1033
$ws.Cells.Item(72, 20).Value = 15

# New row 73: Packham's Triumph / Segunda, 03-12-2021, volumen 200
$ws.Cells.Item(73, 1).Value = 4
$ws.Cells.Item(73, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(73, 3).Value = "Los Lagos"
$ws.Cells.Item(73, 4).Value = 44533
$ws.Cells.Item(73, 5).Value = 10
$ws.Cells.Item(73, 6).Value = "Fruta"
$ws.Cells.Item(73, 7).Value = 100104
$ws.Cells.Item(73, 8).Value = "Frutos de pepita"
$ws.Cells.Item(73, 9).Value = 100104005
$ws.Cells.Item(73, 10).Value = "Pera"
$ws.Cells.Item(73, 11).Value = "Packham's Triumph"
$ws.Cells.Item(73, 12).Value = "Segunda"
$ws.Cells.Item(73, 13).Value = 200
$ws.Cells.Item(73, 14).Value = 14000
$ws.Cells.Item(73, 15).Value = 14000
$ws.Cells.Item(73, 16).Value = 14000
$ws.Cells.Item(73, 17).Value = "`$/caja 15 kilos empedrada"
$ws.Cells.Item(73, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(73, 19).Value = 933
$ws.Cells.Item(73, 20).Value = 15
